$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.764.74"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = "'1.857.83"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').Value = "'1.019"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -1.19%  '
$ws.Range('D5').Value = "'320.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D7').Value = "'0.4371"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').Value = "'0.3777"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').Value = "'0.07424"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = "'0.8835"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').Value = "'21.54"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = "'1.865.82"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').Value = "'6.745"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = "'5.489"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.80%  '
$ws.Range('D15').Value = "'0.07102"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = "'88.18"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.81%  '
$ws.Range('D17').Value = "'1.023"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = "'0.000009028"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = "'1.018"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').Value = "'15.48"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = "'27.766.36"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = "'5.276"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = "'11.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = "'2.088.16"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('D25').Value = "'2.035"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.06%  '
$ws.Range('D26').Value = "'157.06"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').Value = "'18.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('D28').Value = "'5.441"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('D29').Value = "'1.993"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('D30').Value = "'120.84"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.97%  '
$ws.Range('D31').Value = "'0.09043"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').Value = "'1.217"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').Value = "'0.7700"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('D34').Value = "'3.037"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.13%  '
$ws.Range('D35').Value = "'4.557"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = "'1.019"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').Value = "'1.138"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').Value = "'0.01982"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = "'0.05310"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = "'2.870"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = "'0.5179"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = "'6.953"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('D43').Value = "'0.1677"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').Value = "'8.718"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').Value = "'110.13"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = "'10.72"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').Value = "'1.713"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('D48').Value = "'0.4726"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('D49').Value = "'1.020"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = "'0.06467"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('D51').Value = "'1.851"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.27%  '
